# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型"
# sheets to reflect newly generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 645
    $ws.Range("F3").Value = 3862
    $ws.Range("F5").Value = 728
}
